# Outlets Settings.xlsx - "Create supplier should allow only saving of details"
#
# The Reports(existing) row (row 13) is reworded:
#   - E13 now explains that only Owner-permission users get the weekly
#     summary email, with a couple of bold call-outs.
#   - F13's "Actual output" note is expanded to mention that only Owner
#     users receive the weekly emails.
#   - Row 13 grows taller to fit the extra text.
#   - The sheet's scroll position / active selection moves down to the
#     Reports(new) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E13: rebuild the rich text describing the weekly email report ---
$run1 = "Enter the Weekly email report - summary of outlet activity as"
$run2 = "(Sent to all users with Owner permissions)"
$run3 = " and select day "
$run4 = "Send every"
$run5 = " dropdown ,                                                                   *.And Email list of items that are below par and                               *.Email outlet consumption report after stock count to"

$e13 = $ws.Range("E13")
$e13.Value = $run1 + $run2 + $run3 + $run4 + $run5

$pos = 1

$e13.Characters($pos, $run1.Length).Font.Bold = $false
$pos = $pos + $run1.Length

$e13.Characters($pos, $run2.Length).Font.Bold = $true
$pos = $pos + $run2.Length

$e13.Characters($pos, $run3.Length).Font.Bold = $false
$pos = $pos + $run3.Length

$e13.Characters($pos, $run4.Length).Font.Bold = $true
$pos = $pos + $run4.Length

$e13.Characters($pos, $run5.Length).Font.Bold = $false

# --- F13: the actual-output note now calls out that only Owners get the email ---
$ws.Range("F13").Value = "should give email address it get displayed, only Owner users receive the weekly emails."

# --- Row 13 needs more height for the longer text ---
$ws.Rows(13).RowHeight = 96.75

# --- Scroll / selection moves to the Reports(new) row ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("E14").Select()
